$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed shared strings in the exact order needed so that new <si> entries
# land at the same indices as in the target workbook.
$ws.Range("A148").Value = "4.a"
$ws.Range("B153").Value = "Query 10 (see notes)"
$ws.Range("B149").Value = "GEQO effort = 1"
$ws.Range("A153").Value = "4.b"
$ws.Range("A158").Value = "4.c"
$ws.Range("B154").Value = "GEQO effort = 5"
$ws.Range("B159").Value = "GEQO effort = 10"
$ws.Range("B148").Value = "Query 10 (see notes) - 10m tuples"
$ws.Range("B168").Value = "work_mem = 1 Mb"
$ws.Range("B165").Value = "work_mem = 512 kB"
$ws.Range("B173").Value = "work_mem = 4 Mb"

# Fill in the rest of the new rows (values, formulas, styles).
$ws.Range("A148").Value = "4.a"
$ws.Range("A148").Font.Bold = $true
$ws.Range("B148").Value = "Query 10 (see notes) - 10m tuples"
$ws.Range("C148").Value = 1
$ws.Range("D148").Value = 12932

$ws.Range("A149").Value = "4.a"
$ws.Range("A149").Font.Bold = $true
$ws.Range("B149").Value = "GEQO effort = 1"
$ws.Range("C149").Value = 2
$ws.Range("D149").Value = 16544

$ws.Range("A150").Value = "4.a"
$ws.Range("A150").Font.Bold = $true
$ws.Range("C150").Value = 3
$ws.Range("D150").Value = 15971

$ws.Range("A151").Value = "4.a"
$ws.Range("A151").Font.Bold = $true
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 16082

$ws.Range("A152").Value = "4.a"
$ws.Range("A152").Font.Bold = $true
$ws.Range("C152").Value = 5
$ws.Range("D152").Value = 14770
$ws.Range("E152").Formula = "=MAX(D148:D152)"
$ws.Range("F152").Formula = "=MIN(D148:D152)"
$ws.Range("G152").Formula = "=AVERAGE(D148:D152)"
$ws.Range("H152").Formula = "=STDEV.S(D148:D152)"
$ws.Range("I152").Formula = "=((SUM(D148:D152) - E152 - F152) / (COUNT(D148:D152) - 2))"
$ws.Range("E152:I152").Font.Bold = $true

$ws.Range("A153").Value = "4.b"
$ws.Range("A153").Font.Bold = $true
$ws.Range("B153").Value = "Query 10 (see notes)"
$ws.Range("C153").Value = 1
$ws.Range("D153").Value = 12717

$ws.Range("A154").Value = "4.b"
$ws.Range("A154").Font.Bold = $true
$ws.Range("B154").Value = "GEQO effort = 5"
$ws.Range("C154").Value = 2
$ws.Range("D154").Value = 14357

$ws.Range("A155").Value = "4.b"
$ws.Range("A155").Font.Bold = $true
$ws.Range("C155").Value = 3
$ws.Range("D155").Value = 14186

$ws.Range("A156").Value = "4.b"
$ws.Range("A156").Font.Bold = $true
$ws.Range("C156").Value = 4
$ws.Range("D156").Value = 13227

$ws.Range("A157").Value = "4.b"
$ws.Range("A157").Font.Bold = $true
$ws.Range("C157").Value = 5
$ws.Range("D157").Value = 12827
$ws.Range("E157").Formula = "=MAX(D153:D157)"
$ws.Range("F157").Formula = "=MIN(D153:D157)"
$ws.Range("G157").Formula = "=AVERAGE(D153:D157)"
$ws.Range("H157").Formula = "=STDEV.S(D153:D157)"
$ws.Range("I157").Formula = "=((SUM(D153:D157) - E157 - F157) / (COUNT(D153:D157) - 2))"
$ws.Range("E157:I157").Font.Bold = $true

$ws.Range("A158").Value = "4.c"
$ws.Range("A158").Font.Bold = $true
$ws.Range("B158").Value = "Query 10 (see notes)"
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 13094

$ws.Range("A159").Value = "4.c"
$ws.Range("A159").Font.Bold = $true
$ws.Range("B159").Value = "GEQO effort = 10"
$ws.Range("C159").Value = 2
$ws.Range("D159").Value = 13214

$ws.Range("A160").Value = "4.c"
$ws.Range("A160").Font.Bold = $true
$ws.Range("C160").Value = 3
$ws.Range("D160").Value = 11733

$ws.Range("A161").Value = "4.c"
$ws.Range("A161").Font.Bold = $true
$ws.Range("C161").Value = 4
$ws.Range("D161").Value = 11951

$ws.Range("A162").Value = "4.c"
$ws.Range("A162").Font.Bold = $true
$ws.Range("C162").Value = 5
$ws.Range("D162").Value = 11979
$ws.Range("E162").Formula = "=MAX(D158:D162)"
$ws.Range("F162").Formula = "=MIN(D158:D162)"
$ws.Range("G162").Formula = "=AVERAGE(D158:D162)"
$ws.Range("H162").Formula = "=STDEV.S(D158:D162)"
$ws.Range("I162").Formula = "=((SUM(D158:D162) - E162 - F162) / (COUNT(D158:D162) - 2))"
$ws.Range("E162:I162").Font.Bold = $true

$ws.Range("A163").Value = "4.c"
$ws.Range("A163").Font.Bold = $true
$ws.Range("B163").Value = "Query 10 (see notes)"
$ws.Range("C163").Value = 1
$ws.Range("D163").Value = 11334

$ws.Range("A164").Value = "4.c"
$ws.Range("A164").Font.Bold = $true
$ws.Range("B164").Value = "GEQO effort = 10"
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 15318

$ws.Range("A165").Value = "4.c"
$ws.Range("A165").Font.Bold = $true
$ws.Range("B165").Value = "work_mem = 512 kB"
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 14860

$ws.Range("A166").Value = "4.c"
$ws.Range("A166").Font.Bold = $true
$ws.Range("C166").Value = 4
$ws.Range("D166").Value = 14216

$ws.Range("A167").Value = "4.c"
$ws.Range("A167").Font.Bold = $true
$ws.Range("C167").Value = 5
$ws.Range("D167").Value = 15438
$ws.Range("E167").Formula = "=MAX(D163:D167)"
$ws.Range("F167").Formula = "=MIN(D163:D167)"
$ws.Range("G167").Formula = "=AVERAGE(D163:D167)"
$ws.Range("H167").Formula = "=STDEV.S(D163:D167)"
$ws.Range("I167").Formula = "=((SUM(D163:D167) - E167 - F167) / (COUNT(D163:D167) - 2))"
$ws.Range("E167:I167").Font.Bold = $true

$ws.Range("A168").Value = "4.c"
$ws.Range("A168").Font.Bold = $true
$ws.Range("B168").Value = "work_mem = 1 Mb"
$ws.Range("C168").Value = 6
$ws.Range("D168").Value = 11170

$ws.Range("A169").Value = "4.c"
$ws.Range("A169").Font.Bold = $true
$ws.Range("C169").Value = 7
$ws.Range("D169").Value = 13838

$ws.Range("A170").Value = "4.c"
$ws.Range("A170").Font.Bold = $true
$ws.Range("C170").Value = 8
$ws.Range("D170").Value = 11443

$ws.Range("A171").Value = "4.c"
$ws.Range("A171").Font.Bold = $true
$ws.Range("C171").Value = 9
$ws.Range("D171").Value = 11705

$ws.Range("A172").Value = "4.c"
$ws.Range("A172").Font.Bold = $true
$ws.Range("C172").Value = 10
$ws.Range("D172").Value = 11720
$ws.Range("E172").Formula = "=MAX(D168:D172)"
$ws.Range("F172").Formula = "=MIN(D168:D172)"
$ws.Range("G172").Formula = "=AVERAGE(D168:D172)"
$ws.Range("H172").Formula = "=STDEV.S(D168:D172)"
$ws.Range("I172").Formula = "=((SUM(D168:D172) - E172 - F172) / (COUNT(D168:D172) - 2))"
$ws.Range("E172:I172").Font.Bold = $true

$ws.Range("A173").Value = "4.c"
$ws.Range("A173").Font.Bold = $true
$ws.Range("B173").Value = "work_mem = 4 Mb"
$ws.Range("C173").Value = 11
$ws.Range("D173").Value = 10033

$ws.Range("A174").Value = "4.c"
$ws.Range("A174").Font.Bold = $true
$ws.Range("C174").Value = 12
$ws.Range("D174").Value = 12205

$ws.Range("A175").Value = "4.c"
$ws.Range("A175").Font.Bold = $true
$ws.Range("C175").Value = 13
$ws.Range("D175").Value = 12117

$ws.Range("A176").Value = "4.c"
$ws.Range("A176").Font.Bold = $true
$ws.Range("C176").Value = 14
$ws.Range("D176").Value = 12128

$ws.Range("A177").Value = "4.c"
$ws.Range("A177").Font.Bold = $true
$ws.Range("C177").Value = 15
$ws.Range("D177").Value = 12488
$ws.Range("E177").Formula = "=MAX(D173:D177)"
$ws.Range("F177").Formula = "=MIN(D173:D177)"
$ws.Range("G177").Formula = "=AVERAGE(D173:D177)"
$ws.Range("H177").Formula = "=STDEV.S(D173:D177)"
$ws.Range("I177").Formula = "=((SUM(D173:D177) - E177 - F177) / (COUNT(D173:D177) - 2))"
$ws.Range("E177:I177").Font.Bold = $true

# Update selection/scroll to mirror the author leaving off after the new data.
$ws.Range("D178").Select()
$excel.ActiveWindow.ScrollRow = 149
$excel.ActiveWindow.ScrollColumn = 1
